$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 9
$ws.Range("C6").Value = 8
$ws.Range("C7").Value = 5
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 4
$ws.Range("C10").Value = 6
$ws.Range("C12").Value = 8
$ws.Range("C13").Value = 5
$ws.Range("C14").Value = 7
$ws.Range("C17").Value = 3
$ws.Range("C18").Value = 4
